$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.811.53"
$ws.Range("E2").Value = "  -4.03%  "
$ws.Range("D3").Value = "3.432.93"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'577.23"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "'164.59"
$ws.Range("E6").Value = "  -6.27%  "
$ws.Range("E7").Value = "  -4.53%  "
$ws.Range("D8").Value = "3.427.37"
$ws.Range("E8").Value = "  -5.17%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -5.61%  "
$ws.Range("D11").Value = "'6.70"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -8.46%  "
$ws.Range("D13").Value = "'46.09"
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("D14").Value = "'0.0000269"
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").Value = "3.978.19"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").Value = "'609.45"
$ws.Range("E16").Value = "  -9.41%  "
$ws.Range("E17").Value = "  -8.54%  "
$ws.Range("D18").Value = "67.801.50"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").Value = "3.423.66"
$ws.Range("E19").Value = "  -5.43%  "
$ws.Range("D21").Value = "'17.06"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").Value = "'10.86"
$ws.Range("E22").Value = "  -5.02%  "
$ws.Range("D23").Value = "'0.866"
$ws.Range("E23").Value = "  -7.73%  "
$ws.Range("D24").Value = "'15.36"
$ws.Range("E24").Value = "  -9.75%  "
$ws.Range("D25").Value = "'94.76"
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("D26").Value = "'3.70"
$ws.Range("E26").Value = "  -5.21%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'2.56"
$ws.Range("E28").Value = "  -7.86%  "
$ws.Range("D29").Value = "'8.89"
$ws.Range("E29").Value = "  -9.15%  "
$ws.Range("D30").Value = "'31.95"
$ws.Range("E30").Value = "  -7.23%  "
$ws.Range("E31").Value = "  -9.43%  "
$ws.Range("E32").Value = "  -7.91%  "
$ws.Range("D33").Value = "'1.27"
$ws.Range("E33").Value = "  -7.19%  "
$ws.Range("E34").Value = "  -11.11%  "
$ws.Range("D35").Value = "'602.63"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("E36").Value = "  -4.43%  "
$ws.Range("D37").Value = "'56.51"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  -14.25%  "
$ws.Range("D40").Value = "'0.0997"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("D41").Value = "'0.0432"
$ws.Range("E41").Value = "  -4.61%  "
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "3.344.11"
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("E44").Value = "  -7.87%  "
$ws.Range("D45").Value = "'32.21"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("E46").Value = "  -7.06%  "
$ws.Range("E47").Value = "  -7.16%  "
$ws.Range("E48").Value = "  -10.11%  "
$ws.Range("E49").Value = "  -5.66%  "
$ws.Range("D50").Value = "'132.73"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'5.61"
$ws.Range("E51").Value = "  +11.82%  "

# Clear formats on forced-text numeric cells to remove the quote-prefix style marker
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
